$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old summary block (rows 1-2: the "total"/"free" counters) that
# sat above the evaluation table - it's no longer part of the form.
$ws.Range("B1:E2").ClearContents()

# The "nothing" column (E) in the header/table is dropped too.
$ws.Range("E3").ClearContents()

# Reset the view back to the top of the sheet with the next free header
# cell selected (previously it was scrolled down to D66).
$ws.Range("E3").Select()
